# "deleted pics and added readme"
# Tabelle1 is the README-style index sheet: it used to list only 3 stimulus
# files (208, 219, 226). It is replaced with the current, sorted list of
# stimuli actually present after the pruning of images described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$values = @(
    "ImageFile",
    "Stimuli/208.jpg",
    "Stimuli/219.jpg",
    "Stimuli/226.jpg",
    "Stimuli/253.jpg",
    "Stimuli/254.jpg",
    "Stimuli/326.jpg",
    "Stimuli/1301.jpg",
    "Stimuli/3350.jpg",
    "Stimuli/6242.jpg",
    "Stimuli/6410.jpg",
    "Stimuli/6555.jpg",
    "Stimuli/6825.jpg",
    "Stimuli/6940.jpg",
    "Stimuli/8230.jpg",
    "Stimuli/9041.jpg",
    "Stimuli/9140.jpg",
    "Stimuli/9340.jpg",
    "Stimuli/9409.jpg",
    "Stimuli/9570.jpg",
    "Stimuli/9800.jpg"
)

# Clear out anything currently below the new list (sheet previously only had
# 4 rows, so this is a no-op today but keeps the script idempotent/safe).
$ws.Columns.Item(1).ClearContents()

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match the new selection state left behind in the sheet (header row
# excluded, full data block selected).
$ws.Range("A2:A21").Select() | Out-Null
